$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value2 = 0.1991150442477876
$ws.Cells.Item(2, 3).Value2 = 0.5707964601769911
$ws.Cells.Item(2, 10).Value2 = 0.02654867256637168
$ws.Cells.Item(2, 16).Value2 = 0.1371681415929203
$ws.Cells.Item(2, 19).Value2 = 0.06637168141592921
$ws.Cells.Item(3, 2).Value2 = 0.007633587786259542
$ws.Cells.Item(3, 3).Value2 = 0.03053435114503817
$ws.Cells.Item(3, 10).Value2 = 0.01526717557251908
$ws.Cells.Item(3, 16).Value2 = 0.7633587786259542
$ws.Cells.Item(3, 19).Value2 = 0.183206106870229
$ws.Cells.Item(4, 10).Value2 = 0.09677419354838709
$ws.Cells.Item(4, 16).Value2 = 0.6129032258064516
$ws.Cells.Item(4, 19).Value2 = 0.2903225806451613
$ws.Cells.Item(6, 2).Value2 = 0.06896551724137931
$ws.Cells.Item(6, 4).Value2 = 0.01724137931034483
$ws.Cells.Item(6, 5).Value2 = 0.005747126436781609
$ws.Cells.Item(6, 6).Value2 = 0.04022988505747126
$ws.Cells.Item(6, 10).Value2 = 0.2126436781609195
$ws.Cells.Item(6, 15).Value2 = 0.02873563218390805
$ws.Cells.Item(6, 17).Value2 = 0.1379310344827586
$ws.Cells.Item(6, 18).Value2 = 0.09770114942528736
$ws.Cells.Item(6, 19).Value2 = 0.3908045977011494
$ws.Cells.Item(7, 2).Value2 = 0.09271523178807947
$ws.Cells.Item(7, 4).Value2 = 0.01986754966887417
$ws.Cells.Item(7, 6).Value2 = 0.04635761589403974
$ws.Cells.Item(7, 10).Value2 = 0.152317880794702
$ws.Cells.Item(7, 15).Value2 = 0.03311258278145696
$ws.Cells.Item(7, 17).Value2 = 0.1390728476821192
$ws.Cells.Item(7, 18).Value2 = 0.0728476821192053
$ws.Cells.Item(7, 19).Value2 = 0.4437086092715232
$ws.Cells.Item(8, 2).Value2 = 0.07517084282460136
$ws.Cells.Item(8, 4).Value2 = 0.01366742596810934
$ws.Cells.Item(8, 6).Value2 = 0.07289293849658314
$ws.Cells.Item(8, 10).Value2 = 0.132118451025057
$ws.Cells.Item(8, 15).Value2 = 0.02277904328018223
$ws.Cells.Item(8, 17).Value2 = 0.1685649202733485
$ws.Cells.Item(8, 18).Value2 = 0.07517084282460136
$ws.Cells.Item(8, 19).Value2 = 0.4396355353075171
$ws.Cells.Item(9, 2).Value2 = 0.08426966292134831
$ws.Cells.Item(9, 4).Value2 = 0.01123595505617977
$ws.Cells.Item(9, 6).Value2 = 0.06741573033707865
$ws.Cells.Item(9, 10).Value2 = 0.0898876404494382
$ws.Cells.Item(9, 15).Value2 = 0.02247191011235955
$ws.Cells.Item(9, 17).Value2 = 0.1235955056179775
$ws.Cells.Item(9, 18).Value2 = 0.1179775280898876
$ws.Cells.Item(9, 19).Value2 = 0.4831460674157304
$ws.Cells.Item(10, 2).Value2 = 0.0951111111111111
$ws.Cells.Item(10, 4).Value2 = 0.01511111111111111
$ws.Cells.Item(10, 5).Value2 = 0.0008888888888888889
$ws.Cells.Item(10, 6).Value2 = 0.06311111111111112
$ws.Cells.Item(10, 10).Value2 = 0.104
$ws.Cells.Item(10, 15).Value2 = 0.009777777777777778
$ws.Cells.Item(10, 17).Value2 = 0.1982222222222222
$ws.Cells.Item(10, 18).Value2 = 0.08977777777777778
$ws.Cells.Item(10, 19).Value2 = 0.424
$ws.Cells.Item(11, 7).Value2 = 0.1403508771929824
$ws.Cells.Item(11, 10).Value2 = 0.1052631578947368
$ws.Cells.Item(11, 11).Value2 = 0.1798245614035088
$ws.Cells.Item(11, 12).Value2 = 0.5657894736842105
$ws.Cells.Item(11, 19).Value2 = 0.008771929824561403
$ws.Cells.Item(12, 7).Value2 = 0.746268656716418
$ws.Cells.Item(12, 10).Value2 = 0.2014925373134328
$ws.Cells.Item(12, 12).Value2 = 0.04477611940298507
$ws.Cells.Item(12, 19).Value2 = 0.007462686567164179
$ws.Cells.Item(13, 7).Value2 = 0.696969696969697
$ws.Cells.Item(13, 10).Value2 = 0.303030303030303
$ws.Cells.Item(14, 7).Value2 = 0.75
$ws.Cells.Item(14, 10).Value2 = 0.25
$ws.Cells.Item(15, 6).Value2 = 0.01025641025641026
$ws.Cells.Item(15, 8).Value2 = 0.2512820512820513
$ws.Cells.Item(15, 9).Value2 = 0.08205128205128205
$ws.Cells.Item(15, 10).Value2 = 0.3487179487179487
$ws.Cells.Item(15, 11).Value2 = 0.04615384615384616
$ws.Cells.Item(15, 13).Value2 = 0.02051282051282051
$ws.Cells.Item(15, 14).Value2 = 0.005128205128205128
$ws.Cells.Item(15, 15).Value2 = 0.03076923076923077
$ws.Cells.Item(15, 19).Value2 = 0.2051282051282051
$ws.Cells.Item(16, 6).Value2 = 0.006802721088435374
$ws.Cells.Item(16, 8).Value2 = 0.2312925170068027
$ws.Cells.Item(16, 9).Value2 = 0.06802721088435375
$ws.Cells.Item(16, 10).Value2 = 0.4217687074829932
$ws.Cells.Item(16, 11).Value2 = 0.1020408163265306
$ws.Cells.Item(16, 13).Value2 = 0.01360544217687075
$ws.Cells.Item(16, 14).Value2 = 0.006802721088435374
$ws.Cells.Item(16, 15).Value2 = 0.06122448979591837
$ws.Cells.Item(16, 19).Value2 = 0.08843537414965986
$ws.Cells.Item(17, 6).Value2 = 0.005524861878453038
$ws.Cells.Item(17, 8).Value2 = 0.2071823204419889
$ws.Cells.Item(17, 9).Value2 = 0.1132596685082873
$ws.Cells.Item(17, 10).Value2 = 0.4475138121546962
$ws.Cells.Item(17, 11).Value2 = 0.09392265193370165
$ws.Cells.Item(17, 13).Value2 = 0.005524861878453038
$ws.Cells.Item(17, 14).Value2 = 0.002762430939226519
$ws.Cells.Item(17, 15).Value2 = 0.06077348066298342
$ws.Cells.Item(17, 19).Value2 = 0.06353591160220995
$ws.Cells.Item(18, 6).Value2 = 0.01092896174863388
$ws.Cells.Item(18, 8).Value2 = 0.2349726775956284
$ws.Cells.Item(18, 9).Value2 = 0.1038251366120219
$ws.Cells.Item(18, 10).Value2 = 0.4316939890710382
$ws.Cells.Item(18, 11).Value2 = 0.04371584699453552
$ws.Cells.Item(18, 13).Value2 = 0.04918032786885246
$ws.Cells.Item(18, 15).Value2 = 0.07103825136612021
$ws.Cells.Item(18, 19).Value2 = 0.0546448087431694
$ws.Cells.Item(19, 6).Value2 = 0.01176470588235294
$ws.Cells.Item(19, 8).Value2 = 0.2153846153846154
$ws.Cells.Item(19, 9).Value2 = 0.08416289592760182
$ws.Cells.Item(19, 10).Value2 = 0.4072398190045249
$ws.Cells.Item(19, 11).Value2 = 0.1049773755656109
$ws.Cells.Item(19, 13).Value2 = 0.01447963800904977
$ws.Cells.Item(19, 14).Value2 = 0.0009049773755656109
$ws.Cells.Item(19, 15).Value2 = 0.07963800904977375
$ws.Cells.Item(19, 19).Value2 = 0.08144796380090498
